# "edited and uploaded test plan"
# Update the Implementation-phase / Testing-phase contribution notes in
# Table6 (Group Name ... Testing phase ... Maintainance phase, rows 52-57)
# on the Gantt planner sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contribution notes (cell-entry order matters for shared-string layout,
# matching how the author actually typed them in Excel).
$ws.Range("D54").Value = "Game, Table, Lobby, Player, Dealer"
$ws.Range("E55").Value = "Hand, Message, Shoe Test"
$ws.Range("E54").Value = "Table, Lobby, Player, Dealer, Message"
$ws.Range("D55").Value = "Login, Hand, Server, Client, Deck, Card"
$ws.Range("D53").Value = "Server (in the start) , player , Penals, testing"

# Widen columns C and D to fit the longer test-plan text.
$ws.Columns.Item(3).ColumnWidth = 53.75
$ws.Columns.Item(4).ColumnWidth = 46.625

# Leave the selection on D57, and scroll the viewport down a bit, matching
# where the author ended up after editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("D57").Select()
